$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Consumer checking and savings"
$ws.Range("D2").Value = "Business checking and savings"

$ws.Range("A2").Select()
